$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the updated Price cells are plain decimal-looking strings
# (e.g. '1.001', '18.60', '1.000') where the trailing digits are meaningful
# text, not numeric precision. Excel's General number format would silently
# coerce these into numbers on assignment (dropping the trailing zero, e.g.
# '1.000' -> 1), so those specific cells are temporarily switched to the
# Text number format before the value is written, then restored to the
# workbook's normal/default style so no visible formatting changes stick.
$textForceCells = @(
    'D4',
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D13',
    'D14',
    'D15',
    'D16',
    'D18',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D50',
    'D51'
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated coin values (price + 1h volume change) row by row,
# matching the latest pull from the coinranking feed.

# Row 2
$ws.Range('D2').Value = '27.083.47'
$ws.Range('E2').Value = '  -0.42%  '

# Row 3
$ws.Range('D3').Value = '1.825.49'
$ws.Range('E3').Value = '  +0.44%  '

# Row 4
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.46%  '

# Row 5
$ws.Range('D5').Value = '312.43'
$ws.Range('E5').Value = '  -0.40%  '

# Row 6
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.37%  '

# Row 7
$ws.Range('D7').Value = '0.4565'
$ws.Range('E7').Value = '  +7.20%  '

# Row 8
$ws.Range('D8').Value = '0.3732'
$ws.Range('E8').Value = '  +1.93%  '

# Row 9
$ws.Range('D9').Value = '0.07327'
$ws.Range('E9').Value = '  +1.87%  '

# Row 10
$ws.Range('D10').Value = '0.8596'

# Row 11
$ws.Range('D11').Value = '20.96'
$ws.Range('E11').Value = '  +0.13%  '

# Row 12
$ws.Range('D12').Value = '1.788.58'
$ws.Range('E12').Value = '  -2.03%  '

# Row 13
$ws.Range('D13').Value = '6.692'
$ws.Range('E13').Value = '  +0.75%  '

# Row 14
$ws.Range('D14').Value = '92.94'
$ws.Range('E14').Value = '  +5.92%  '

# Row 15
$ws.Range('D15').Value = '5.337'
$ws.Range('E15').Value = '  +0.78%  '

# Row 16
$ws.Range('D16').Value = '0.07085'
$ws.Range('E16').Value = '  +0.07%  '

# Row 17
$ws.Range('E17').Value = '  -0.41%  '

# Row 18
$ws.Range('D18').Value = '0.000008826'
$ws.Range('E18').Value = '  -0.30%  '

# Row 19
$ws.Range('E19').Value = '  -0.34%  '

# Row 20
$ws.Range('E20').Value = '  +0.03%  '

# Row 21
$ws.Range('D21').Value = '27.129.82'
$ws.Range('E21').Value = '  -0.36%  '

# Row 22
$ws.Range('E22').Value = '  +1.60%  '

# Row 23
$ws.Range('E23').Value = '  +1.48%  '

# Row 24
$ws.Range('D24').Value = '1.997'
$ws.Range('E24').Value = '  -0.18%  '

# Row 25
$ws.Range('D25').Value = '151.69'
$ws.Range('E25').Value = '  -0.75%  '

# Row 26
$ws.Range('D26').Value = '2.225'
$ws.Range('E26').Value = '  +5.51%  '

# Row 27
$ws.Range('D27').Value = '18.60'
$ws.Range('E27').Value = '  +1.91%  '

# Row 28
$ws.Range('D28').Value = '5.261'
$ws.Range('E28').Value = '  +0.79%  '

# Row 29
$ws.Range('D29').Value = '117.50'
$ws.Range('E29').Value = '  +1.25%  '

# Row 30
$ws.Range('D30').Value = '0.08856'
$ws.Range('E30').Value = '  -0.19%  '

# Row 31
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '0.7598'
$ws.Range('E31').Value = '  +0.23%  '

# Row 32
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.191'
$ws.Range('E32').Value = '  -0.01%  '

# Row 33
$ws.Range('D33').Value = '2.955'
$ws.Range('E33').Value = '  +4.61%  '

# Row 34
$ws.Range('D34').Value = '4.470'
$ws.Range('E34').Value = '  +0.34%  '

# Row 35
$ws.Range('E35').Value = '  -0.35%  '

# Row 36
$ws.Range('D36').Value = '1.103'
$ws.Range('E36').Value = '  -0.81%  '

# Row 37
$ws.Range('D37').Value = '0.01965'
$ws.Range('E37').Value = '  +0.65%  '

# Row 38
$ws.Range('D38').Value = '0.05279'
$ws.Range('E38').Value = '  +0.76%  '

# Row 39
$ws.Range('D39').Value = '0.5359'
$ws.Range('E39').Value = '  +7.24%  '

# Row 40
$ws.Range('D40').Value = '7.176'
$ws.Range('E40').Value = '  +1.97%  '

# Row 41
$ws.Range('D41').Value = '2.890'
$ws.Range('E41').Value = '  -0.25%  '

# Row 42
$ws.Range('D42').Value = '0.1708'
$ws.Range('E42').Value = '  +2.05%  '

# Row 43
$ws.Range('D43').Value = '0.5184'
$ws.Range('E43').Value = '  +10.63%  '

# Row 44
$ws.Range('D44').Value = '8.622'
$ws.Range('E44').Value = '  +0.38%  '

# Row 45
$ws.Range('D45').Value = '10.68'
$ws.Range('E45').Value = '  +1.72%  '

# Row 46
$ws.Range('D46').Value = '1.974'
$ws.Range('E46').Value = '  +9.66%  '

# Row 47
$ws.Range('D47').Value = '105.95'
$ws.Range('E47').Value = '  -0.21%  '

# Row 48
$ws.Range('D48').Value = '1.674'
$ws.Range('E48').Value = '  +1.31%  '

# Row 49
$ws.Range('E49').Value = '  -0.36%  '

# Row 50
$ws.Range('D50').Value = '0.06355'
$ws.Range('E50').Value = '  -0.70%  '

# Row 51
$ws.Range('D51').Value = '0.9230'
$ws.Range('E51').Value = '  +1.49%  '

# Restore the default (unstyled) look on the cells we temporarily
# formatted as Text above, so only the values themselves changed.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}

Write-Output "Updated $($textForceCells.Count) numeric-text cells and applied all row changes."
